$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 6, pushing the existing
# rows 6, 7, 8 (price history) down to rows 8, 9, 10. Excel carries the
# row/cell formatting (e.g. the date style on column D) down with the
# shifted rows, and copies it up into the freshly inserted rows too.
$ws.Rows.Item(6).Resize(2).Insert()

# Fill in the new weekly price entries for 2021-09-28 (serial 44467) in
# the two newly inserted rows, following the same layout as the rows
# beneath them.

# Row 6: "Primera" quality
$ws.Cells.Item(6, 1).Value = 11
$ws.Cells.Item(6, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(6, 3).Value = "Bíobío"
$ws.Cells.Item(6, 4).Value = 44467
$ws.Cells.Item(6, 5).Value = 8
$ws.Cells.Item(6, 6).Value = "Fruta"
$ws.Cells.Item(6, 7).Value = 100107
$ws.Cells.Item(6, 8).Value = "Otros"
$ws.Cells.Item(6, 9).Value = 100107002
$ws.Cells.Item(6, 10).Value = "Chirimoya"
$ws.Cells.Item(6, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(6, 12).Value = "Primera"
$ws.Cells.Item(6, 13).Value = 100
$ws.Cells.Item(6, 14).Value = 2700
$ws.Cells.Item(6, 15).Value = 2800
$ws.Cells.Item(6, 16).Value = 2750
$ws.Cells.Item(6, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(6, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(6, 19).Value = 2750
$ws.Cells.Item(6, 20).Value = 1

# Row 7: "Segunda" quality
$ws.Cells.Item(7, 1).Value = 11
$ws.Cells.Item(7, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(7, 3).Value = "Bíobío"
$ws.Cells.Item(7, 4).Value = 44467
$ws.Cells.Item(7, 5).Value = 8
$ws.Cells.Item(7, 6).Value = "Fruta"
$ws.Cells.Item(7, 7).Value = 100107
$ws.Cells.Item(7, 8).Value = "Otros"
$ws.Cells.Item(7, 9).Value = 100107002
$ws.Cells.Item(7, 10).Value = "Chirimoya"
$ws.Cells.Item(7, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(7, 12).Value = "Segunda"
$ws.Cells.Item(7, 13).Value = 50
$ws.Cells.Item(7, 14).Value = 2500
$ws.Cells.Item(7, 15).Value = 2500
$ws.Cells.Item(7, 16).Value = 2500
$ws.Cells.Item(7, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(7, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(7, 19).Value = 2500
$ws.Cells.Item(7, 20).Value = 1

$ws.Range("D6:D7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
